# Insert a new data row at row 107 (a new weekly price observation),
# pushing the previously-existing rows 107:164 down to 108:165.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(107).Insert()

# Populate the newly inserted row with the new observation's data.
$ws.Range('A107').Value = 9
$ws.Range('B107').Value = 'Vega Central Mapocho de Santiago'
$ws.Range('C107').Value = 'Metropolitana'
$ws.Range('D107').Value = 44518
$ws.Range('E107').Value = 13
$ws.Range('F107').Value = 100112030
$ws.Range('G107').Value = 'Poroto granado'
$ws.Range('H107').Value = 'Sin especificar'
$ws.Range('I107').Value = 'Primera'
$ws.Range('J107').Value = 34
$ws.Range('K107').Value = 35000
$ws.Range('L107').Value = 37000
$ws.Range('M107').Value = 36000
$ws.Range('N107').Value = '$/malla 25 kilos'
$ws.Range('O107').Value = 'Provincia de Limarí'
$ws.Range('P107').Value = 1440
$ws.Range('Q107').Value = 25
$ws.Range('R107').Value = 'Hortaliza'
